$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("renta_total_pv_pbi")

# Insert a new column before R (renta_pv shifts to S, renta_pbi -> T,
# renta_usd_tcc -> U, renta_usd_tcp -> V), making room for the new
# "proporcion_subsidios" metric.
$ws.Range("R1").EntireColumn.Insert()

# Header for the newly inserted column, matching the bold, centered
# header style used by the rest of row 1.
$ws.Range("R1").Value = "proporcion_subsidios"
$ws.Range("R1").Font.Bold = $true
$ws.Range("R1").HorizontalAlignment = -4108

# Fill proporcion_subsidios = subsidios / (subsidios + renta_total)
# for every data row that has both a "subsidios" (J) and "renta_total"
# (Q) figure.
for ($row = 53; $row -le 109; $row++) {
    $j = $ws.Cells.Item($row, 10).Value2
    $q = $ws.Cells.Item($row, 17).Value2
    $ws.Cells.Item($row, 18).Value = $j / ($j + $q)
}
